# Update the two-digit division problems in the single table on the page.
# The table has 20 rows x 5 columns; only every 4th row (1, 5, 9, 13, 17 in
# 1-based Word indexing) actually holds a "NN÷N=" expression, the others are
# blank working-space rows. We target each populated cell directly by its
# (row, column) table coordinate rather than a document-wide Find/Replace,
# since several of the new values reuse text that is also an *old* value
# elsewhere in the table (e.g. "74÷7=" is both a replacement target and a
# replacement result), which would make a sequential Find/Replace pass
# ambiguous/order-dependent.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Each entry: row, column, expected current text, new text.
$edits = @(
    @(1, 1, "19÷4=", "54÷7="),
    @(1, 2, "23÷5=", "74÷7="),
    @(1, 3, "96÷8=", "22÷2="),
    @(1, 4, "83÷6=", "79÷3="),
    @(1, 5, "35÷7=", "83÷7="),

    @(5, 1, "66÷5=", "65÷2="),
    @(5, 2, "44÷9=", "91÷5="),
    @(5, 3, "94÷6=", "96÷7="),
    @(5, 4, "82÷6=", "85÷3="),
    @(5, 5, "87÷9=", "67÷9="),

    @(9, 1, "11÷4=", "25÷2="),
    @(9, 2, "38÷9=", "56÷2="),
    @(9, 3, "33÷4=", "25÷4="),
    @(9, 4, "23÷9=", "71÷4="),
    @(9, 5, "12÷9=", "58÷7="),

    @(13, 1, "44÷7=", "54÷3="),
    @(13, 2, "81÷7=", "41÷3="),
    @(13, 3, "74÷7=", "80÷5="),
    @(13, 4, "55÷8=", "62÷6="),
    @(13, 5, "32÷9=", "97÷5="),

    @(17, 1, "37÷7=", "69÷9="),
    @(17, 2, "32÷5=", "72÷4="),
    @(17, 3, "50÷3=", "76÷6="),
    @(17, 4, "38÷5=", "89÷9="),
    @(17, 5, "84÷3=", "14÷8=")
)

$mismatches = New-Object System.Collections.ArrayList
$updated = 0

foreach ($edit in $edits) {
    $row = $edit[0]
    $col = $edit[1]
    $old = $edit[2]
    $new = $edit[3]

    $cell = $t.Cell($row, $col)
    $r = $cell.Range
    # Cell.Range.Text includes a trailing paragraph mark + cell mark;
    # trim those off so we only compare/replace the visible content.
    $r.MoveEnd(1, -1) | Out-Null

    if ($r.Text -ne $old) {
        [void]$mismatches.Add("Row ${row} Col ${col}: expected '$old' but found '" + $r.Text + "'")
    } else {
        $r.Text = $new
        $updated = $updated + 1
    }
}

if ($mismatches.Count -gt 0) {
    Write-Output "Updated $updated of $($edits.Count) cells; mismatches encountered:"
    foreach ($m in $mismatches) { Write-Output $m }
} else {
    Write-Output "Updated all $updated cells successfully."
}
